$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.294.39'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '2.089.06'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '342.73'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Value = '0.5224'
$ws.Range('E7').Value = '  +1.68%  '
$ws.Range('D8').Value = '0.4401'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '54.36'
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('D10').Value = '0.09348'
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').Value = '1.167'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = '24.76'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '8.620'
$ws.Range('E13').Value = '  +3.99%  '
$ws.Range('D14').Value = '6.897'
$ws.Range('E14').Value = '  +2.05%  '
$ws.Range('D15').Value = '2.046.83'
$ws.Range('E15').Value = '  -2.49%  '
$ws.Range('D16').Value = '101.29'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '0.00001156'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '1.004'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = '21.10'
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').Value = '0.06668'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').Value = '6.326'
$ws.Range('E21').Value = '  +2.07%  '
$ws.Range('D23').Value = '30.301.29'
$ws.Range('E23').Value = '  +1.84%  '
$ws.Range('D24').Value = '12.50'
$ws.Range('E24').Value = '  -1.00%  '
$ws.Range('D25').Value = '2.298'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '21.79'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').Value = '162.13'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = '2.509'
$ws.Range('E28').Value = '  -0.89%  '
$ws.Range('D29').Value = '133.00'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '1.130'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('B31').Value = 'ARBITRUM'
$ws.Range('C31').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D31').Value = '1.662'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.1047'
$ws.Range('E32').Value = '  -0.28%  '
$ws.Range('D33').Value = '6.218'
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').Value = '6.655'
$ws.Range('E34').Value = '  +10.47%  '
$ws.Range('D35').Value = '3.867'
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('E36').Value = '  -2.64%  '
$ws.Range('D37').Value = '0.02629'
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +3.64%  '
$ws.Range('D41').Value = '12.50'
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('D42').Value = '0.2206'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('D43').Value = '0.6796'
$ws.Range('E43').Value = '  +2.06%  '
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '2.331'
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '1.370'
$ws.Range('E47').Value = '  +18.19%  '
$ws.Range('D48').Value = '3.623'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '0.00000000347'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = '1.213'
$ws.Range('E50').Value = '  +8.49%  '
$ws.Range('E51').Value = '  -0.52%  '

Write-Output "Applied 93 cell updates"
